$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Gnai2"
$ws.Cells.Item(2, 3).Value = "Adcy1"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 195.0792385
$ws.Cells.Item(2, 8).Value = 390.158477
$ws.Cells.Item(2, 9).Value = 0.2640605522989327
$ws.Cells.Item(2, 10).Value = 0.1982306263353075
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.007396666666666667
$ws.Cells.Item(2, 14).Value = 0.02219
$ws.Cells.Item(2, 15).Value = 0.03674254717418517
$ws.Cells.Item(2, 16).Value = 0.05411957523816028
$ws.Cells.Item(2, 17).Value = 1.442936100771667
$ws.Cells.Item(2, 18).Value = 8.65761660463
$ws.Cells.Item(2, 19).Value = 0.009702257299684925
$ws.Cells.Item(2, 20).Value = 0.01072815729646131

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Gnai2"
$ws.Cells.Item(3, 3).Value = "Adcy1"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 195.0792385
$ws.Cells.Item(3, 8).Value = 390.158477
$ws.Cells.Item(3, 9).Value = 0.2640605522989327
$ws.Cells.Item(3, 10).Value = 0.1982306263353075
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.193914
$ws.Cells.Item(3, 14).Value = 0.387828
$ws.Cells.Item(3, 15).Value = 0.9632574528258149
$ws.Cells.Item(3, 16).Value = 0.9458804247618398
$ws.Cells.Item(3, 17).Value = 37.828595454489
$ws.Cells.Item(3, 18).Value = 151.314381817956
$ws.Cells.Item(3, 19).Value = 0.2543582949992477
$ws.Cells.Item(3, 20).Value = 0.1875024690388462

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Gnai2"
$ws.Cells.Item(4, 3).Value = "Adcy1"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 62.40792233333334
$ws.Cells.Item(4, 8).Value = 187.223767
$ws.Cells.Item(4, 9).Value = 0.08447577797556809
$ws.Cells.Item(4, 10).Value = 0.09512412720758515
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.007396666666666667
$ws.Cells.Item(4, 14).Value = 0.02219
$ws.Cells.Item(4, 15).Value = 0.03674254717418517
$ws.Cells.Item(4, 16).Value = 0.05411957523816028
$ws.Cells.Item(4, 17).Value = 0.4616105988588889
$ws.Cells.Item(4, 18).Value = 4.15449538973
$ws.Cells.Item(4, 19).Value = 0.003103855257343303
$ws.Cells.Item(4, 20).Value = 0.005148077359375234

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Gnai2"
$ws.Cells.Item(5, 3).Value = "Adcy1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 62.40792233333334
$ws.Cells.Item(5, 8).Value = 187.223767
$ws.Cells.Item(5, 9).Value = 0.08447577797556809
$ws.Cells.Item(5, 10).Value = 0.09512412720758515
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.193914
$ws.Cells.Item(5, 14).Value = 0.387828
$ws.Cells.Item(5, 15).Value = 0.9632574528258149
$ws.Cells.Item(5, 16).Value = 0.9458804247618398
$ws.Cells.Item(5, 17).Value = 12.101769851346
$ws.Cells.Item(5, 18).Value = 72.61061910807601
$ws.Cells.Item(5, 19).Value = 0.08137192271822478
$ws.Cells.Item(5, 20).Value = 0.08997604984820992

# Row 6
$ws.Cells.Item(6, 1).Value = "M1"
$ws.Cells.Item(6, 2).Value = "Gnai2"
$ws.Cells.Item(6, 3).Value = "Adcy1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 142.8621113333333
$ws.Cells.Item(6, 8).Value = 428.586334
$ws.Cells.Item(6, 9).Value = 0.1933791023142199
$ws.Cells.Item(6, 10).Value = 0.2177549443006804
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.007396666666666667
$ws.Cells.Item(6, 14).Value = 0.02219
$ws.Cells.Item(6, 15).Value = 0.03674254717418517
$ws.Cells.Item(6, 16).Value = 0.05411957523816028
$ws.Cells.Item(6, 17).Value = 1.056703416828889
$ws.Cells.Item(6, 18).Value = 9.51033075146
$ws.Cells.Item(6, 19).Value = 0.007105240789281806
$ws.Cells.Item(6, 20).Value = 0.01178480509156207

# Row 7
$ws.Cells.Item(7, 1).Value = "M1"
$ws.Cells.Item(7, 2).Value = "Gnai2"
$ws.Cells.Item(7, 3).Value = "Adcy1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 142.8621113333333
$ws.Cells.Item(7, 8).Value = 428.586334
$ws.Cells.Item(7, 9).Value = 0.1933791023142199
$ws.Cells.Item(7, 10).Value = 0.2177549443006804
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.193914
$ws.Cells.Item(7, 14).Value = 0.387828
$ws.Cells.Item(7, 15).Value = 0.9632574528258149
$ws.Cells.Item(7, 16).Value = 0.9458804247618398
$ws.Cells.Item(7, 17).Value = 27.702963457092
$ws.Cells.Item(7, 18).Value = 166.217780742552
$ws.Cells.Item(7, 19).Value = 0.1862738615249381
$ws.Cells.Item(7, 20).Value = 0.2059701392091183

# Row 8
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Gnai2"
$ws.Cells.Item(8, 3).Value = "Adcy1"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 141.6168416666667
$ws.Cells.Item(8, 8).Value = 424.850525
$ws.Cells.Item(8, 9).Value = 0.1916934970264942
$ws.Cells.Item(8, 10).Value = 0.2158568649262854
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.007396666666666667
$ws.Cells.Item(8, 14).Value = 0.02219
$ws.Cells.Item(8, 15).Value = 0.03674254717418517
$ws.Cells.Item(8, 16).Value = 0.05411957523816028
$ws.Cells.Item(8, 17).Value = 1.047492572194444
$ws.Cells.Item(8, 18).Value = 9.427433149750001
$ws.Cells.Item(8, 19).Value = 0.007043307357480487
$ws.Cells.Item(8, 20).Value = 0.0116820818420515

# Row 9
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Gnai2"
$ws.Cells.Item(9, 3).Value = "Adcy1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 141.6168416666667
$ws.Cells.Item(9, 8).Value = 424.850525
$ws.Cells.Item(9, 9).Value = 0.1916934970264942
$ws.Cells.Item(9, 10).Value = 0.2158568649262854
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.193914
$ws.Cells.Item(9, 14).Value = 0.387828
$ws.Cells.Item(9, 15).Value = 0.9632574528258149
$ws.Cells.Item(9, 16).Value = 0.9458804247618398
$ws.Cells.Item(9, 17).Value = 27.46148823495
$ws.Cells.Item(9, 18).Value = 164.7689294097
$ws.Cells.Item(9, 19).Value = 0.1846501896690137
$ws.Cells.Item(9, 20).Value = 0.2041747830842339

# Row 10
$ws.Cells.Item(10, 1).Value = "Neutro"
$ws.Cells.Item(10, 2).Value = "Gnai2"
$ws.Cells.Item(10, 3).Value = "Adcy1"
$ws.Cells.Item(10, 4).Value = "FAPs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 143.783834
$ws.Cells.Item(10, 8).Value = 431.351502
$ws.Cells.Item(10, 9).Value = 0.1946267522348261
$ws.Cells.Item(10, 10).Value = 0.2191598631141254
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.007396666666666667
$ws.Cells.Item(10, 14).Value = 0.02219
$ws.Cells.Item(10, 15).Value = 0.03674254717418517
$ws.Cells.Item(10, 16).Value = 0.05411957523816028
$ws.Cells.Item(10, 17).Value = 1.063521092153333
$ws.Cells.Item(10, 18).Value = 9.57168982938
$ws.Cells.Item(10, 19).Value = 0.007151082625346548
$ws.Cells.Item(10, 20).Value = 0.01186083870098982

# Row 11
$ws.Cells.Item(11, 1).Value = "Neutro"
$ws.Cells.Item(11, 2).Value = "Gnai2"
$ws.Cells.Item(11, 3).Value = "Adcy1"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 143.783834
$ws.Cells.Item(11, 8).Value = 431.351502
$ws.Cells.Item(11, 9).Value = 0.1946267522348261
$ws.Cells.Item(11, 10).Value = 0.2191598631141254
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.193914
$ws.Cells.Item(11, 14).Value = 0.387828
$ws.Cells.Item(11, 15).Value = 0.9632574528258149
$ws.Cells.Item(11, 16).Value = 0.9458804247618398
$ws.Cells.Item(11, 17).Value = 27.881698386276
$ws.Cells.Item(11, 18).Value = 167.290190317656
$ws.Cells.Item(11, 19).Value = 0.1874756696094796
$ws.Cells.Item(11, 20).Value = 0.2072990244131356

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Gnai2"
$ws.Cells.Item(12, 3).Value = "Adcy1"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 53.01711450000001
$ws.Cells.Item(12, 8).Value = 106.034229
$ws.Cells.Item(12, 9).Value = 0.07176431814995911
$ws.Cells.Item(12, 10).Value = 0.05387357411601602
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.007396666666666667
$ws.Cells.Item(12, 14).Value = 0.02219
$ws.Cells.Item(12, 15).Value = 0.03674254717418517
$ws.Cells.Item(12, 16).Value = 0.05411957523816028
$ws.Cells.Item(12, 17).Value = 0.3921499235850001
$ws.Cells.Item(12, 18).Value = 2.35289954151
$ws.Cells.Item(12, 19).Value = 0.002636803845048106
$ws.Cells.Item(12, 20).Value = 0.002915614947720333

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Gnai2"
$ws.Cells.Item(13, 3).Value = "Adcy1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 53.01711450000001
$ws.Cells.Item(13, 8).Value = 106.034229
$ws.Cells.Item(13, 9).Value = 0.07176431814995911
$ws.Cells.Item(13, 10).Value = 0.05387357411601602
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.193914
$ws.Cells.Item(13, 14).Value = 0.387828
$ws.Cells.Item(13, 15).Value = 0.9632574528258149
$ws.Cells.Item(13, 16).Value = 0.9458804247618398
$ws.Cells.Item(13, 17).Value = 10.280760741153
$ws.Cells.Item(13, 18).Value = 41.123042964612
$ws.Cells.Item(13, 19).Value = 0.069127514304911
$ws.Cells.Item(13, 20).Value = 0.05095795916829569
